$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 11.73582396971049
$ws.Range("C2").Value = 8.198143628860791
$ws.Range("D2").Value = 5.742490069362098
$ws.Range("E2").Value = 11.71729909709096
$ws.Range("F2").Value = 27.40552900906796
$ws.Range("I2").Value = 24.69008135183571
$ws.Range("K2").Value = 9.102014730398244
$ws.Range("L2").Value = 9.757681587473193
$ws.Range("O2").Value = 24.62061181298066
$ws.Range("B3").Value = 11.44897535467203
$ws.Range("C3").Value = 8.169600342180633
$ws.Range("D3").Value = 5.693058072245429
$ws.Range("E3").Value = 11.73032061351807
$ws.Range("F3").Value = 27.46051221529218
$ws.Range("I3").Value = 24.78850513011865
$ws.Range("K3").Value = 8.895974339975837
$ws.Range("L3").Value = 9.742494873551044
$ws.Range("O3").Value = 24.70489543553953
$ws.Range("B4").Value = 11.27087887610324
$ws.Range("C4").Value = 8.15195975560099
$ws.Range("D4").Value = 5.662036351539615
$ws.Range("E4").Value = 11.74045278583219
$ws.Range("F4").Value = 27.50119698236893
$ws.Range("I4").Value = 24.85356120080091
$ws.Range("K4").Value = 8.768145275271308
$ws.Range("L4").Value = 9.734843055466042
$ws.Range("O4").Value = 24.76173537559975
$ws.Range("B5").Value = 11.19791374239899
$ws.Range("C5").Value = 8.144744233902701
$ws.Range("D5").Value = 5.64923129031952
$ws.Range("E5").Value = 11.74511918157899
$ws.Range("F5").Value = 27.51951383276393
$ws.Range("I5").Value = 24.88123389307195
$ws.Range("K5").Value = 8.715794003208661
$ws.Range("L5").Value = 9.732148133322758
$ws.Range("O5").Value = 24.78617573377882
$ws.Range("B6").Value = 11.18577761330029
$ws.Range("C6").Value = 8.1435445405838
$ws.Range("D6").Value = 5.647095297006338
$ws.Range("E6").Value = 11.7459264943064
$ws.Range("F6").Value = 27.52266014856468
$ws.Range("I6").Value = 24.8858990704265
$ws.Range("K6").Value = 8.707087573488245
$ws.Range("L6").Value = 9.73172627735755
$ws.Range("O6").Value = 24.790311137311
$ws.Range("B7").Value = 11.26989627568582
$ws.Range("C7").Value = 8.151862550527115
$ws.Range("D7").Value = 5.661864313086894
$ws.Range("E7").Value = 11.74051354230508
$ws.Range("F7").Value = 27.50143698031799
$ws.Range("I7").Value = 24.85392970130947
$ws.Range("K7").Value = 8.767440202934933
$ws.Range("L7").Value = 9.734804993839127
$ws.Range("O7").Value = 24.76205981669801
$ws.Range("B8").Value = 11.63738960603847
$ws.Range("C8").Value = 8.188326981283295
$ws.Range("D8").Value = 5.725588060955352
$ws.Range("E8").Value = 11.72134550686427
$ws.Range("F8").Value = 27.42304781744336
$ws.Range("I8").Value = 24.7230577915754
$ws.Range("K8").Value = 9.031288079031155
$ws.Range("L8").Value = 9.752099641275667
$ws.Range("O8").Value = 24.64861520939244
$ws.Range("B9").Value = 12.33801069803126
$ws.Range("C9").Value = 8.258832972700109
$ws.Range("D9").Value = 5.845013524147199
$ws.Range("E9").Value = 11.70070291672773
$ws.Range("F9").Value = 27.32442202873888
$ws.Range("I9").Value = 24.50313892265602
$ws.Range("K9").Value = 9.535221206957978
$ws.Range("L9").Value = 9.799162481231694
$ws.Range("O9").Value = 24.46662766829249
$ws.Range("B10").Value = 12.83505702715018
$ws.Range("C10").Value = 8.309906016622556
$ws.Range("D10").Value = 5.929095737029099
$ws.Range("E10").Value = 11.69584752845379
$ws.Range("F10").Value = 27.28572048647957
$ws.Range("I10").Value = 24.36399854537263
$ws.Range("K10").Value = 9.893472513480644
$ws.Range("L10").Value = 9.841570157991839
$ws.Range("O10").Value = 24.35772293067243
$ws.Range("B11").Value = 13.0562948324768
$ws.Range("C11").Value = 8.332957245999264
$ws.Range("D11").Value = 5.966493746375496
$ws.Range("E11").Value = 11.69587015974644
$ws.Range("F11").Value = 27.27546862142072
$ws.Range("I11").Value = 24.30558206826069
$ws.Range("K11").Value = 10.05312446611081
$ws.Range("L11").Value = 9.862517659500204
$ws.Range("O11").Value = 24.31359021063855
$ws.Range("B12").Value = 13.13929349360031
$ws.Range("C12").Value = 8.341657745441267
$ws.Range("D12").Value = 5.980528123935162
$ws.Range("E12").Value = 11.69619866660733
$ws.Range("F12").Value = 27.27264478405326
$ws.Range("I12").Value = 24.28416390671112
$ws.Range("K12").Value = 10.11304854139369
$ws.Range("L12").Value = 9.870683780162366
$ws.Range("O12").Value = 24.29765803813193
$ws.Range("B13").Value = 13.12145409337995
$ws.Range("C13").Value = 8.339785241601049
$ws.Range("D13").Value = 5.977511315825139
$ws.Range("E13").Value = 11.69611370416224
$ws.Range("F13").Value = 27.2732058722737
$ws.Range("I13").Value = 24.28874540940989
$ws.Range("K13").Value = 10.1001673489313
$ws.Range("L13").Value = 9.86891473615265
$ws.Range("O13").Value = 24.30105459524027
$ws.Range("B14").Value = 13.06313922928507
$ws.Range("C14").Value = 8.333673625762202
$ws.Range("D14").Value = 5.967650945332069
$ws.Range("E14").Value = 11.69589078013829
$ws.Range("F14").Value = 27.27521508973587
$ws.Range("I14").Value = 24.30380588892429
$ws.Range("K14").Value = 10.05806543410756
$ws.Range("L14").Value = 9.863184832366704
$ws.Range("O14").Value = 24.31226381810283
$ws.Range("B15").Value = 13.02731597378005
$ws.Range("C15").Value = 8.329926304354846
$ws.Range("D15").Value = 5.961594441137859
$ws.Range("E15").Value = 11.69579586738433
$ws.Range("F15").Value = 27.27658363228073
$ws.Range("I15").Value = 24.31312244250928
$ws.Range("K15").Value = 10.03220585237898
$ws.Range("L15").Value = 9.85970540830551
$ws.Range("O15").Value = 24.31923142382128
$ws.Range("B16").Value = 12.82049371165012
$ws.Range("C16").Value = 8.308395702083878
$ws.Range("D16").Value = 5.926634126053562
$ws.Range("E16").Value = 11.69589087462449
$ws.Range("F16").Value = 27.28653853487032
$ws.Range("I16").Value = 24.36791448939996
$ws.Range("K16").Value = 9.882967272160116
$ws.Range("L16").Value = 9.840234144951786
$ws.Range("O16").Value = 24.36071613118327
$ws.Range("B17").Value = 12.69230921763506
$ws.Range("C17").Value = 8.295139081521965
$ws.Range("D17").Value = 5.904965440039645
$ws.Range("E17").Value = 11.69652005782227
$ws.Range("F17").Value = 27.29452983418522
$ws.Range("I17").Value = 24.40277821004348
$ws.Range("K17").Value = 9.790523539726527
$ws.Range("L17").Value = 9.828710275913631
$ws.Range("O17").Value = 24.38755250920038
$ws.Range("B18").Value = 12.61812670288055
$ws.Range("C18").Value = 8.28749721729886
$ws.Range("D18").Value = 5.892422338246242
$ws.Range("E18").Value = 11.69709204121114
$ws.Range("F18").Value = 27.29981837048076
$ws.Range("I18").Value = 24.4232900745291
$ws.Range("K18").Value = 9.737043064677101
$ws.Range("L18").Value = 9.822238204709723
$ws.Range("O18").Value = 24.40349704181806
$ws.Range("B19").Value = 12.5929342791677
$ws.Range("C19").Value = 8.284906962644287
$ws.Range("D19").Value = 5.888161900185502
$ws.Range("E19").Value = 11.69732181435818
$ws.Range("F19").Value = 27.30172781794409
$ws.Range("I19").Value = 24.43031386457893
$ws.Range("K19").Value = 9.718884132868849
$ws.Range("L19").Value = 9.82007382114389
$ws.Range("O19").Value = 24.40898293250893
$ws.Range("B20").Value = 12.70600227370196
$ws.Range("C20").Value = 8.296552048233838
$ws.Range("D20").Value = 5.907280409116843
$ws.Range("E20").Value = 11.69643134248035
$ws.Range("F20").Value = 27.29360750303324
$ws.Range("I20").Value = 24.39901937545085
$ws.Range("K20").Value = 9.800396757698698
$ws.Range("L20").Value = 9.829920881023337
$ws.Range("O20").Value = 24.38464303959164
$ws.Range("B21").Value = 13.08028946478013
$ws.Range("C21").Value = 8.335469547405024
$ws.Range("D21").Value = 5.970550673059632
$ws.Range("E21").Value = 11.6959475836091
$ws.Range("F21").Value = 27.27459620778699
$ws.Range("I21").Value = 24.29936317458354
$ws.Range("K21").Value = 10.07044666477054
$ws.Range("L21").Value = 9.86486153614454
$ws.Range("O21").Value = 24.30895021364379
$ws.Range("B22").Value = 13.32033338800913
$ws.Range("C22").Value = 8.360737183977793
$ws.Range("D22").Value = 6.011156412622016
$ws.Range("E22").Value = 11.69749570257581
$ws.Range("F22").Value = 27.26833984885126
$ws.Range("I22").Value = 24.23832966848157
$ws.Range("K22").Value = 10.24381234995752
$ws.Range("L22").Value = 9.889057568427631
$ws.Range("O22").Value = 24.26402754985232
$ws.Range("B23").Value = 13.19266038037031
$ws.Range("C23").Value = 8.347267433235508
$ws.Range("D23").Value = 5.989554158180232
$ws.Range("E23").Value = 11.69649921249005
$ws.Range("F23").Value = 27.27111444672276
$ws.Range("I23").Value = 24.27052904732971
$ws.Range("K23").Value = 10.15158725248625
$ws.Range("L23").Value = 9.876020738770436
$ws.Range("O23").Value = 24.28758687792556
$ws.Range("B24").Value = 12.69981315557995
$ws.Range("C24").Value = 8.295913309532267
$ws.Range("D24").Value = 5.906234077481113
$ws.Range("E24").Value = 11.69647079569442
$ws.Range("F24").Value = 27.29402232653114
$ws.Range("I24").Value = 24.40071728609058
$ws.Range("K24").Value = 9.79593410987982
$ws.Range("L24").Value = 9.829373089069026
$ws.Range("O24").Value = 24.38595680382798
$ws.Range("B25").Value = 12.15120890833902
$ws.Range("C25").Value = 8.239878446319697
$ws.Range("D25").Value = 5.813328195577565
$ws.Range("E25").Value = 11.70447453801877
$ws.Range("F25").Value = 27.34518447720552
$ws.Range("I25").Value = 24.55869670598284
$ws.Range("K25").Value = 9.400740579616683
$ws.Range("L25").Value = 9.785041076537901
$ws.Range("O25").Value = 24.51151349514871
